# Applies the commit: two betting-odds rows (2 and 10) had their match
# data swapped back to the correct match, and two new match rows (101, 102)
# were appended at the end of the sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2: replace odds/teams/url (columns F:V) ---
$ws.Range("F2").Value = "URSL Vise"
$ws.Range("G2").Value = 0
$ws.Range("H2").Value = "Charleroi B"
$ws.Range("I2").Value = 1
$ws.Range("J2").Value = 1.74
$ws.Range("K2").Value = "30/08/2023 09:12"
$ws.Range("L2").Value = 1.93
$ws.Range("M2").Value = "30/08/2023 19:18"
$ws.Range("N2").Value = 3.73
$ws.Range("O2").Value = "30/08/2023 09:12"
$ws.Range("P2").Value = 3.35
$ws.Range("Q2").Value = "30/08/2023 19:18"
$ws.Range("R2").Value = 4.27
$ws.Range("S2").Value = "30/08/2023 09:12"
$ws.Range("T2").Value = 3.9
$ws.Range("U2").Value = "30/08/2023 19:18"
$ws.Range("V2").Value = "https://www.betexplorer.com/football/belgium/national-division-1/ursl-vise-charleroi/bPegWb68/"

# --- Row 10: replace odds/teams/url (columns F:V) ---
$ws.Range("F10").Value = "Tienen"
$ws.Range("G10").Value = 1
$ws.Range("H10").Value = "Heist"
$ws.Range("I10").Value = 1
$ws.Range("J10").Value = 2.92
$ws.Range("K10").Value = "30/08/2023 09:12"
$ws.Range("L10").Value = 2.53
$ws.Range("M10").Value = "30/08/2023 19:24"
$ws.Range("N10").Value = 3.36
$ws.Range("O10").Value = "30/08/2023 09:12"
$ws.Range("P10").Value = 3.55
$ws.Range("Q10").Value = "30/08/2023 18:06"
$ws.Range("R10").Value = 2.23
$ws.Range("S10").Value = "30/08/2023 09:12"
$ws.Range("T10").Value = 2.54
$ws.Range("U10").Value = "30/08/2023 19:24"
$ws.Range("V10").Value = "https://www.betexplorer.com/football/belgium/national-division-1/tienen-heist/CSacVILE/"

# --- New row 101 (values first, then copy formatting from row 100) ---
$ws.Range("A101").Value = 100
$ws.Range("B101").Value = "belgium"
$ws.Range("C101").Value = "national-division-1"
$ws.Range("D101").Value = "2023-2024"
$ws.Range("E101").Value = 45242.60416666666
$ws.Range("F101").Value = "Antwerp B"
$ws.Range("G101").Value = 2
$ws.Range("H101").Value = "Saint Eloois"
$ws.Range("I101").Value = 1
$ws.Range("J101").Value = 2.45
$ws.Range("K101").Value = "11/11/2023 02:43"
$ws.Range("L101").Value = 2.89
$ws.Range("M101").Value = "12/11/2023 14:17"
$ws.Range("N101").Value = 3.12
$ws.Range("O101").Value = "11/11/2023 02:43"
$ws.Range("P101").Value = 3.32
$ws.Range("Q101").Value = "12/11/2023 14:17"
$ws.Range("R101").Value = 2.59
$ws.Range("S101").Value = "11/11/2023 02:43"
$ws.Range("T101").Value = 2.36
$ws.Range("U101").Value = "12/11/2023 14:17"
$ws.Range("V101").Value = "https://www.betexplorer.com/football/belgium/national-division-1/antwerp-saint-eloois/zoexN9oo/"
$ws.Range("A100:V100").Copy()
$ws.Range("A101:V101").PasteSpecial(-4122)

# --- New row 102 (values first, then copy formatting from row 100) ---
$ws.Range("A102").Value = 101
$ws.Range("B102").Value = "belgium"
$ws.Range("C102").Value = "national-division-1"
$ws.Range("D102").Value = "2023-2024"
$ws.Range("E102").Value = 45242.60416666666
$ws.Range("F102").Value = "Namur"
$ws.Range("G102").Value = 1
$ws.Range("H102").Value = "Thes Sport"
$ws.Range("I102").Value = 0
$ws.Range("J102").Value = 2.5
$ws.Range("K102").Value = "11/11/2023 02:43"
$ws.Range("L102").Value = 2.65
$ws.Range("M102").Value = "12/11/2023 14:20"
$ws.Range("N102").Value = 3.2
$ws.Range("O102").Value = "11/11/2023 02:43"
$ws.Range("P102").Value = 3.13
$ws.Range("Q102").Value = "12/11/2023 14:20"
$ws.Range("R102").Value = 2.49
$ws.Range("S102").Value = "11/11/2023 02:43"
$ws.Range("T102").Value = 2.67
$ws.Range("U102").Value = "12/11/2023 14:20"
$ws.Range("V102").Value = "https://www.betexplorer.com/football/belgium/national-division-1/namur-thes-sport/djtXri84/"
$ws.Range("A100:V100").Copy()
$ws.Range("A102:V102").PasteSpecial(-4122)

$excel.CutCopyMode = 0

Write-Output "Applied odds swap on rows 2/10 and appended rows 101-102."
